$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new reading ("2026/01/16", 金, 8) was recorded for the 2026/01/16 date,
# pushing everything from the old row 648 ("2026/12/29 ...") down by one row
# (old 648..689 become new 649..690).
$ws.Rows.Item(648).Insert()

# Column A holds the date as plain text (not a real Excel date), exactly
# like row 647 right above it (which is also "2026/01/16"). Typing the
# string directly into the cell would make Excel auto-convert it into a
# date serial, so instead copy the existing text cell from row 647 - this
# preserves the literal text value/type without touching number formats.
$ws.Cells.Item(647, 1).Copy()
$ws.Cells.Item(648, 1).PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Cells.Item(648, 2).Value = "金"
$ws.Cells.Item(648, 3).Value = 8
$ws.Cells.Item(648, 4).Value = 201
